# M01 Froze TE+Decoder 1
# Update the per-epoch accuracy figures in column B (rows 2-109) of
# "Epoch Accuracy" Sheet1 with the re-run values, and restore the
# normal top-left selection (A2:B109) instead of the whole-sheet
# selection left over from the previous run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.984375,
    0.984375,
    0.9375,
    0.921875,
    0.921875,
    0.875,
    0.890625,
    0.890625,
    0.890625,
    0.859375,
    0.90625,
    0.875,
    0.859375,
    0.921875,
    0.890625,
    0.875,
    0.84375,
    0.859375,
    0.8125,
    0.765625,
    0.8125,
    0.828125,
    0.828125,
    0.84375,
    0.84375,
    0.84375,
    0.84375,
    0.828125,
    0.828125,
    0.8125,
    0.828125,
    0.828125,
    0.84375,
    0.84375,
    0.84375,
    0.84375,
    0.84375,
    0.84375,
    0.84375,
    0.828125,
    0.796875,
    0.796875,
    0.8125,
    0.84375,
    0.84375,
    0.828125,
    0.828125,
    0.859375,
    0.84375,
    0.8125,
    0.765625,
    0.8125,
    0.8125,
    0.765625,
    0.75,
    0.765625,
    0.78125,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.734375,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.71875,
    0.75,
    0.8125,
    0.8125,
    0.734375,
    0.671875,
    0.546875,
    0.8235294117647058
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}

$ws.Range("A2:B109").Select()